$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '63.003.97'

$ws.Range("D3").Value = '2.587.11'
$ws.Range("E3").Value = '  +1.69%  '

$ws.Range("E4").Value = '  -0.03%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '583.19'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +1.59%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '147.50'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +1.06%  '

$ws.Range("E7").Value = '  -0.01%  '

$ws.Range("E8").Value = '  +2.67%  '

$ws.Range("E9").Value = '  +2.70%  '

$ws.Range("E10").Value = '  +2.87%  '

$ws.Range("E11").Value = '  -0.04%  '

$ws.Range("E12").Value = '  -0.05%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '27.36'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  +0.23%  '

$ws.Range("D14").Value = '3.050.69'
$ws.Range("E14").Value = '  +1.76%  '

$ws.Range("D15").Value = '62.860.44'
$ws.Range("E15").Value = '  -0.03%  '

$ws.Range("E16").Value = '  +3.12%  '

$ws.Range("D17").Value = '2.585.56'
$ws.Range("E17").Value = '  +2.33%  '

$ws.Range("E18").Value = '  -0.05%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '343.43'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +2.21%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '4.40'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +1.71%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '6.69'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -0.75%  '

$ws.Range("E22").Value = '  -0.06%  '

$ws.Range("E23").Value = '  -1.18%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '67.35'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +3.32%  '

$ws.Range("D25").Value = '2.720.08'
$ws.Range("E25").Value = '  +2.10%  '

$ws.Range("E26").Value = '  -0.67%  '

$ws.Range("E27").Value = '  +0.69%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '0.998'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  -0.32%  '

$ws.Range("B29").Value = 'Aptos'
$ws.Range("C29").Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '7.87'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  +8.71%  '

$ws.Range("B30").Value = 'InternetComputer(DFINITY)'
$ws.Range("C30").Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '8.34'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +0.40%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '1.44'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -0.76%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '1.94'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +4.60%  '

$ws.Range("E33").Value = '  +1.87%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '467.69'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  +16.80%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '175.03'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  -1.57%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '1.61'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +4.80%  '

$ws.Range("E37").Value = '  +0.06%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.404'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +1.55%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '19.15'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  +0.40%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '4.56'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +5.42%  '

$ws.Range("E42").Value = '  -1.38%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '158.23'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  +5.12%  '

$ws.Range("E44").Value = '  +1.39%  '

$ws.Range("E45").Value = '  +6.56%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '21.32'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +2.84%  '

$ws.Range("E47").Value = '  +1.98%  '

$ws.Range("E48").Value = '  +0.87%  '

$ws.Range("E49").Value = '  -0.16%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '18.45'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +1.58%  '

$ws.Range("B51").Value = 'dogwifhat'
$ws.Range("C51").Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '1.72'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +1.47%  '
